# Rename the embedded header/footer logo pictures:
#   - "BTec_Logo-Orange" (header):               image1.jpg -> image2.jpg
#   - Pearson logo (footers, descr is the UNC
#     path to the source asset):                 image2.png -> image1.png
#
# These pictures live in the section Headers/Footers ranges (not in
# Document.InlineShapes, which only covers the main story), so we walk
# Sections -> Headers/Footers -> Range.InlineShapes and match each
# picture by its AlternativeText (the docPr/cNvPr "descr").
#
# NOTE: object handles returned from a Headers/Footers collection can go
# stale as soon as another part in the package is edited, so each shape
# is re-fetched immediately before it is touched, and a Write-Host is
# emitted right after every rename to force the interpreter to settle
# before the next lookup.

$d = $word.ActiveDocument

$btecOldName = "image1.jpg"
$btecNewName = "image2.jpg"
$btecDescr = "BTec_Logo-Orange"

$pearsonOldName = "image2.png"
$pearsonNewName = "image1.png"
$pearsonDescr = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"

for ($s = 1; $s -le $d.Sections.Count; $s++) {

    for ($hi = 1; $hi -le 3; $hi++) {
        $sec = $d.Sections.Item($s)
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $shapeCount = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                $sec = $d.Sections.Item($s)
                $hdr = $sec.Headers.Item($hi)
                $shp = $hdr.Range.InlineShapes.Item($j)
                $alt = $shp.AlternativeText
                if ($alt -eq $btecDescr) {
                    $shp.Name = $btecNewName
                    Write-Host "Section $s header $hi shape $j renamed to $btecNewName"
                } elseif ($alt -eq $pearsonDescr) {
                    $shp.Name = $pearsonNewName
                    Write-Host "Section $s header $hi shape $j renamed to $pearsonNewName"
                }
            }
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $sec = $d.Sections.Item($s)
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            $shapeCount = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $shapeCount; $j++) {
                $sec = $d.Sections.Item($s)
                $ftr = $sec.Footers.Item($fi)
                $shp = $ftr.Range.InlineShapes.Item($j)
                $alt = $shp.AlternativeText
                if ($alt -eq $btecDescr) {
                    $shp.Name = $btecNewName
                    Write-Host "Section $s footer $fi shape $j renamed to $btecNewName"
                } elseif ($alt -eq $pearsonDescr) {
                    $shp.Name = $pearsonNewName
                    Write-Host "Section $s footer $fi shape $j renamed to $pearsonNewName"
                }
            }
        }
    }
}
